$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = " > INTRODUZIONE ALLA PROGRAMMAZIONE E LABORATORIO [H3, 1A]`n > INTRODUZIONE ALL INTELLIGENZA ARTIFICIALE [H2 Bis, 5C]"
$ws.Range("C2").Value = " > INTRODUZIONE ALLA PROGRAMMAZIONE E LABORATORIO [H3, 1A]`n > INTRODUZIONE ALL INTELLIGENZA ARTIFICIALE [H2 Bis, 5C]"
$ws.Range("D2").Value = " > INTRODUZIONE ALLA PROGRAMMAZIONE E LABORATORIO [H3, 1A]"
$ws.Range("E2").Value = " > NESSUNA LEZIONE"
$ws.Range("F2").Value = " > PROGRAMMAZIONE AVANZATA E PARALLELA [H3, 1B]"
$ws.Range("G2").Value = " > ALGEBRA LINEARE ED ELEMENTI DI GEOMETRIA [H3, 1A]`n > PROGRAMMAZIONE AVANZATA E PARALLELA [H3, 1B]"
$ws.Range("H2").Value = " > ALGEBRA LINEARE ED ELEMENTI DI GEOMETRIA [H3, 1A]`n > ALGORITMI E STRUTTURE DATI [H2 Bis, 3B]`n > INTRODUZIONE AL MACHINE LEARNING [H3, 1B]"
$ws.Range("I2").Value = " > ALGORITMI E STRUTTURE DATI [H2 Bis, 3B]`n > INTRODUZIONE AL MACHINE LEARNING [H3, 1B]"
$ws.Range("B3").Value = " > ANALISI 1 [H3, 1A]`n > PROGRAMMAZIONE AVANZATA E PARALLELA [H3, 2A]"
$ws.Range("C3").Value = " > ANALISI 1 [H3, 1A]`n > PROGRAMMAZIONE AVANZATA E PARALLELA [H3, 2A]"
$ws.Range("D3").Value = " > ARCHITETTURE DEGLI ELABORATORI [H3, 2B]"
$ws.Range("E3").Value = " > ARCHITETTURE DEGLI ELABORATORI [H3, 2B]"
$ws.Range("F3").Value = " > ARCHITETTURE DEGLI ELABORATORI [H3, 2B]"
$ws.Range("G3").Value = " > INTRODUZIONE ALLA PROGRAMMAZIONE E LABORATORIO [H3, 2B]"
$ws.Range("H3").Value = " > INTRODUZIONE ALLA PROGRAMMAZIONE E LABORATORIO [H3, 2B]`n > ALGORITMI E STRUTTURE DATI [H3, 2A]"
$ws.Range("I3").Value = " > INTRODUZIONE ALLA PROGRAMMAZIONE E LABORATORIO [H3, 2B]`n > ALGORITMI E STRUTTURE DATI [H3, 2A]"
$ws.Range("B4").Value = " > SISTEMI COMPLESSI [H3, 1B]"
$ws.Range("C4").Value = " > INFERENZA STATISTICA [H3, 1C]`n > SISTEMI COMPLESSI [H3, 1B]"
$ws.Range("D4").Value = " > INFERENZA STATISTICA [H3, 1C]`n > SISTEMI COMPLESSI [H3, 1B]"
$ws.Range("E4").Value = " > PROGRAMMAZIONE AVANZATA E PARALLELA [H3, 2A]"
$ws.Range("F4").Value = " > PROGRAMMAZIONE AVANZATA E PARALLELA [H3, 2A]"
$ws.Range("G4").Value = " > NESSUNA LEZIONE"
$ws.Range("H4").Value = " > INTRODUZIONE ALLA PROGRAMMAZIONE E LABORATORIO [H3, 2B]"
$ws.Range("I4").Value = " > INTRODUZIONE ALLA PROGRAMMAZIONE E LABORATORIO [H3, 2B]"
$ws.Range("B5").Value = " > ALGORITMI E STRUTTURE DATI [H2 Bis, 2A MORIN]`n > SISTEMI COMPLESSI [H2 Bis, 3B]"
$ws.Range("C5").Value = " > ALGORITMI E STRUTTURE DATI [H2 Bis, 2A MORIN]`n > SISTEMI COMPLESSI [H2 Bis, 3B]"
$ws.Range("D5").Value = " > ANALISI 1 [H3, 1B]`n > INFERENZA STATISTICA [H2 Bis, 2A MORIN]"
$ws.Range("E5").Value = " > ANALISI 1 [H3, 1B]`n > INFERENZA STATISTICA [H2 Bis, 2A MORIN]"
$ws.Range("F5").Value = " > ALGEBRA LINEARE ED ELEMENTI DI GEOMETRIA [H2 Bis, 2A MORIN]`n > INTRODUZIONE AL MACHINE LEARNING [H2 Bis, 5A]"
$ws.Range("G5").Value = " > ALGEBRA LINEARE ED ELEMENTI DI GEOMETRIA [H2 Bis, 2A MORIN]`n > INTRODUZIONE AL MACHINE LEARNING [H2 Bis, 5A]"
$ws.Range("H5").Value = " > ANALISI 1 [H3, 2B]`n > INTRODUZIONE ALL INTELLIGENZA ARTIFICIALE [H3, 1B]"
$ws.Range("I5").Value = " > ANALISI 1 [H3, 2B]`n > INTRODUZIONE ALL INTELLIGENZA ARTIFICIALE [H3, 1B]"
$ws.Range("B6").Value = " > COMPUTABILITà, COMPLESSITà E LOGICA [H3, 2C]`n > INTRODUZIONE ALL INTELLIGENZA ARTIFICIALE [H2 Bis, 5B]"
$ws.Range("C6").Value = " > INTRODUZIONE ALLA PROGRAMMAZIONE E LABORATORIO [H3, 1B]`n > COMPUTABILITà, COMPLESSITà E LOGICA [H3, 2C]`n > INTRODUZIONE ALL INTELLIGENZA ARTIFICIALE [H2 Bis, 5B]"
$ws.Range("D6").Value = " > INTRODUZIONE ALLA PROGRAMMAZIONE E LABORATORIO [H3, 1B]`n > COMPUTABILITà, COMPLESSITà E LOGICA [H3, 2C]`n > INTRODUZIONE ALL INTELLIGENZA ARTIFICIALE [H2 Bis, 5B]"
$ws.Range("E6").Value = " > INTRODUZIONE ALLA PROGRAMMAZIONE E LABORATORIO [H3, 1B]`n > COMPUTABILITà, COMPLESSITà E LOGICA [H3, 2C]"
$ws.Range("F6").Value = " > INTRODUZIONE ALLA PROGRAMMAZIONE E LABORATORIO [H3, 1B]`n > COMPUTABILITà, COMPLESSITà E LOGICA [H3, 2C]"
$ws.Range("G6").Value = " > COMPUTABILITà, COMPLESSITà E LOGICA [H3, 2C]"
$ws.Range("H6").Value = " > ALGEBRA LINEARE ED ELEMENTI DI GEOMETRIA [H3, 1A]`n > INFERENZA STATISTICA [H3, 1C]`n > INTRODUZIONE AL MACHINE LEARNING [H3, 2C]"
$ws.Range("I6").Value = " > ALGEBRA LINEARE ED ELEMENTI DI GEOMETRIA [H3, 1A]`n > INFERENZA STATISTICA [H3, 1C]`n > INTRODUZIONE AL MACHINE LEARNING [H3, 2C]"
